# Updates metadata on the "attribute" sheet:
#  - row 12 and row 13 no longer wrap text in column G (shrinks their
#    custom row heights accordingly)
#  - selection moves to C15

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attribute")

# Turn off word-wrap for the "number of fish" / "days" cells in column G
# of rows 12-13; this also shrinks the (previously auto/custom) row
# heights now that the text no longer needs to wrap.
$ws.Range("G12").WrapText = $false
$ws.Range("G13").WrapText = $false

$ws.Rows.Item(12).RowHeight = 16
$ws.Rows.Item(13).RowHeight = 14

# Move the active selection to C15
$ws.Range("C15").Select()
